$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename the "_old"/"_new" suffixed header labels (row 1) to the new
#    "_FV2310"/"_FV2404" suffixes used by the regenerated merged AHB file.
# ---------------------------------------------------------------------
$ws.Cells.Replace("_old", "_FV2310")
$ws.Cells.Replace("_new", "_FV2404")

# ---------------------------------------------------------------------
# 2. Freeze the header row (row 1) so it stays visible while scrolling.
# ---------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------
# 3. Turn the data range into a proper Excel Table ("Table1") spanning
#    the whole used range, keeping the header row's existing formatting
#    (bold / shaded / bordered / centered+wrapped) intact instead of the
#    banded "TableStyleMedium2" look Excel applies by default.
# ---------------------------------------------------------------------
$hdr = $ws.Range("A1:U1")

# stash a copy of the header formatting (identical across A1:U1) in a
# scratch cell well outside the used range so it survives ClearFormats
$scratch = $ws.Range("W1")
$ws.Range("A1").Copy($scratch)

$hdr.ClearFormats()
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U69"), [System.Type]::Missing, 1)
$lo.Name = "Table1"

# restore the original header look from the stashed copy, then remove it
$scratch.Copy()
$hdr.PasteSpecial(-4122)
$scratch.Clear()

# the workbook ships with no particular table style applied
$lo.TableStyle = ""
